# Rename file and delete note relevant for testing:
# The author's change clears out the helper/testing "id" column (column C),
# which held a simple running index and its "id" header, while leaving the
# rest of the table (columns D:I) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole column C (as the author did before deleting its contents)
# and clear its contents without shifting the remaining columns.
$col = $ws.Range("C1:C1048576")
$col.Select()
$col.ClearContents()
